# Apply the two changes described by the commit:
#  1. Slide 5's table switches from the custom "Table_0" style to the
#     built-in table style {1FF322EB-4C8A-445E-A1A4-A765E1298222}.
#  2. The deck's theme (colour scheme) is swapped from the colourful
#     "Integral / Red Violet" palette to the plain default
#     "Office Theme / Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(5)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{1FF322EB-4C8A-445E-A1A4-A765E1298222}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
function HexToVbRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Slot order exposed by ThemeColorScheme.Colors(1..12):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = HexToVbRgb $officeColors[$i - 1]
}

$notesColorScheme = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesColorScheme.Colors($i).RGB = HexToVbRgb $officeColors[$i - 1]
}
